# Add bookmarks around the "banner plan" and "crosstab" expected-text
# cells so they match the reality the author observed (commit message:
# "updated banner expected and crosstab expected to match reality
# better"). Word renumbers all w:id values by document order whenever
# a bookmark is added, so inserting these three new bookmarks also
# naturally bumps the pre-existing OLE_LINK1 bookmark's id from 0 to 2.

$d = $word.ActiveDocument

function Add-ParagraphBookmark($searchText, $bookmarkName) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    $para = $rng.Paragraphs(1).Range
    # Exclude the trailing paragraph mark so the bookmark wraps only the
    # run content, not the pilcrow.
    $para.MoveEnd(1, -1) | Out-Null
    $d.Bookmarks.Add($bookmarkName, $para) | Out-Null
}

# "S2=3 OR S2=4 OR S2=5" cell -> OLE_LINK2
Add-ParagraphBookmark "S2=3" "OLE_LINK2"

# "HCP Role" header cell -> OLE_LINK3
Add-ParagraphBookmark "Role" "OLE_LINK3"

# "S9=1 OR S9=2 OR S9=3 OR S9=4 OR S9=5" cell -> OLE_LINK4
Add-ParagraphBookmark "S9=1" "OLE_LINK4"
